$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.337264536586076
$ws.Range("C2").Value = 0.2495539691077795
$ws.Range("D2").Value = 0.1111683606676266
$ws.Range("F2").Value = 2.364975423259352
$ws.Range("G2").Value = 0.002527410921510674
$ws.Range("L2").Value = 0.2117761014173212
$ws.Range("M2").Value = 0.4293080002094669
$ws.Range("N2").Value = 2.256116028307801
$ws.Range("B3").Value = 2.210620735654459
$ws.Range("C3").Value = 0.2170596959110185
$ws.Range("D3").Value = 0.1116371452467853
$ws.Range("F3").Value = 2.313137542642536
$ws.Range("G3").Value = 0.002533449005866058
$ws.Range("L3").Value = 0.2103934094145288
$ws.Range("M3").Value = 0.4124278570964819
$ws.Range("N3").Value = 2.266791252660397
$ws.Range("B4").Value = 2.134362324837184
$ws.Range("C4").Value = 0.1970960426278339
$ws.Range("D4").Value = 0.1119670057879461
$ws.Range("F4").Value = 2.282753959586074
$ws.Range("G4").Value = 0.002537349794861247
$ws.Range("L4").Value = 0.2096418277288876
$ws.Range("M4").Value = 0.4023346689815099
$ws.Range("N4").Value = 2.274030355306806
$ws.Range("B5").Value = 2.103661279415689
$ws.Range("C5").Value = 0.1889568277343585
$ws.Range("D5").Value = 0.1121119486657634
$ws.Range("F5").Value = 2.270733354847053
$ws.Range("G5").Value = 0.002538988200766456
$ws.Range("L5").Value = 0.2093600275653102
$ws.Range("M5").Value = 0.3982894909721466
$ws.Range("N5").Value = 2.277151763214093
$ws.Range("B6").Value = 2.098585959099069
$ws.Range("C6").Value = 0.1876050596591199
$ws.Range("D6").Value = 0.1121366505069794
$ws.Range("F6").Value = 2.268759072856014
$ws.Range("G6").Value = 0.00253926320912452
$ws.Range("L6").Value = 0.2093147126921764
$ws.Range("M6").Value = 0.3976218821530466
$ws.Range("N6").Value = 2.277680403252774
$ws.Range("B7").Value = 2.133946765415033
$ws.Range("C7").Value = 0.1969862910749782
$ws.Range("D7").Value = 0.1119689179941652
$ws.Range("F7").Value = 2.282590387201267
$ws.Range("G7").Value = 0.002537371693353791
$ws.Range("L7").Value = 0.2096379281917606
$ws.Range("M7").Value = 0.4022798399018512
$ws.Range("N7").Value = 2.274071758569107
$ws.Range("B8").Value = 2.293284479049134
$ws.Range("C8").Value = 0.238351926834099
$ws.Range("D8").Value = 0.1113212477112242
$ws.Range("F8").Value = 2.346800304825038
$ws.Range("G8").Value = 0.002529452822260357
$ws.Range("L8").Value = 0.2112791237254399
$ws.Range("M8").Value = 0.4234312083639651
$ws.Range("N8").Value = 2.259654359686465
$ws.Range("B9").Value = 2.617798737017608
$ws.Range("C9").Value = 0.3194160369768895
$ws.Range("D9").Value = 0.1103867210829677
$ws.Range("F9").Value = 2.48430136029657
$ws.Range("G9").Value = 0.002515450333093741
$ws.Range("L9").Value = 0.2152714883575868
$ws.Range("M9").Value = 0.4670797934956852
$ws.Range("N9").Value = 2.236844790518461
$ws.Range("B10").Value = 2.863791852350573
$ws.Range("C10").Value = 0.3790046136435876
$ws.Range("D10").Value = 0.1099076773497885
$ws.Range("F10").Value = 2.592566302782245
$ws.Range("G10").Value = 0.002506082016836031
$ws.Range("L10").Value = 0.2186788516577991
$ws.Range("M10").Value = 0.5005013663077449
$ws.Range("N10").Value = 2.223461598558401
$ws.Range("B11").Value = 2.977395882083044
$ws.Range("C11").Value = 0.4061338973652937
$ws.Range("D11").Value = 0.1097354805116737
$ws.Range("F11").Value = 2.643431837776205
$ws.Range("G11").Value = 0.002502017347139329
$ws.Range("L11").Value = 0.2203324825262172
$ws.Range("M11").Value = 0.5160062576869393
$ws.Range("N11").Value = 2.218115915600308
$ws.Range("B12").Value = 3.020662971093657
$ws.Range("C12").Value = 0.4164112883722169
$ws.Range("D12").Value = 0.1096769042582579
$ws.Range("F12").Value = 2.66292862961518
$ws.Range("G12").Value = 0.002500506310045789
$ws.Range("L12").Value = 0.220973604496038
$ws.Range("M12").Value = 0.5219213573524328
$ws.Range("N12").Value = 2.216199226934052
$ws.Range("B13").Value = 3.011333561091362
$ws.Range("C13").Value = 0.4141976700271357
$ws.Range("D13").Value = 0.1096892238846365
$ws.Range("F13").Value = 2.658719143931222
$ws.Range("G13").Value = 0.002500830488996464
$ws.Range("L13").Value = 0.2208348629984727
$ws.Range("M13").Value = 0.5206454830226406
$ws.Range("N13").Value = 2.216607220682306
$ws.Range("B14").Value = 2.980950508644696
$ws.Range("C14").Value = 0.4069793365193277
$ws.Range("D14").Value = 0.1097305281886918
$ws.Range("F14").Value = 2.645031122010039
$ws.Range("G14").Value = 0.002501892469731264
$ws.Range("L14").Value = 0.2203849286644726
$ws.Range("M14").Value = 0.5164920175953256
$ws.Range("N14").Value = 2.217956066989103
$ws.Range("B15").Value = 2.962372366152636
$ws.Range("C15").Value = 0.402558463419382
$ws.Range("D15").Value = 0.1097566934965357
$ws.Range("F15").Value = 2.636677517815968
$ws.Range("G15").Value = 0.002502546626250498
$ws.Range("L15").Value = 0.2201112758492911
$ws.Range("M15").Value = 0.5139536079996461
$ws.Range("N15").Value = 2.218796313122013
$ws.Range("B16").Value = 2.85640205052249
$ws.Range("C16").Value = 0.3772321763704554
$ws.Range("D16").Value = 0.1099198556440513
$ws.Range("F16").Value = 2.589274865827605
$ws.Range("G16").Value = 0.002506351602704848
$ws.Range("L16").Value = 0.2185728702729932
$ws.Range("M16").Value = 0.4994941823042325
$ws.Range("N16").Value = 2.223825962797264
$ws.Range("B17").Value = 2.791830476117013
$ws.Range("C17").Value = 0.3617016518317087
$ws.Range("D17").Value = 0.1100317005664593
$ws.Range("F17").Value = 2.560610592119502
$ws.Range("G17").Value = 0.002508736176679408
$ws.Range("L17").Value = 0.2176556644031251
$ws.Range("M17").Value = 0.4907012419800196
$ws.Range("N17").Value = 2.227102236485351
$ws.Range("B18").Value = 2.754850481854305
$ws.Range("C18").Value = 0.3527709320284202
$ws.Range("D18").Value = 0.1101003310140882
$ws.Range("F18").Value = 2.544275523989853
$ws.Range("G18").Value = 0.002510126275280136
$ws.Range("L18").Value = 0.2171378616897499
$ws.Range("M18").Value = 0.4856721049678754
$ws.Range("N18").Value = 2.229056488977236
$ws.Range("B19").Value = 2.742357053039086
$ws.Range("C19").Value = 0.3497474636484412
$ws.Range("D19").Value = 0.1101243049756278
$ws.Range("F19").Value = 2.538770747472313
$ws.Range("G19").Value = 0.002510600130657382
$ws.Range("L19").Value = 0.2169642162548371
$ws.Range("M19").Value = 0.4839741787815797
$ws.Range("N19").Value = 2.229730132360274
$ws.Range("B20").Value = 2.798687662004738
$ws.Range("C20").Value = 0.3633546861467494
$ws.Range("D20").Value = 0.1100193490266079
$ws.Range("F20").Value = 2.563646215199441
$ws.Range("G20").Value = 0.002508480415685964
$ws.Range("L20").Value = 0.2177522932526728
$ws.Range("M20").Value = 0.4916343305038424
$ws.Range("N20").Value = 2.226746238667346
$ws.Range("B21").Value = 2.989867999242847
$ws.Range("C21").Value = 0.4090994166933797
$ws.Range("D21").Value = 0.1097182156985355
$ws.Range("F21").Value = 2.649045221924013
$ws.Range("G21").Value = 0.00250157977700042
$ws.Range("L21").Value = 0.2205166799021825
$ws.Range("M21").Value = 0.5177108007959959
$ws.Range("N21").Value = 2.217556950917924
$ws.Range("B22").Value = 3.116260761721492
$ws.Range("C22").Value = 0.4390207027225301
$ws.Range("D22").Value = 0.109560084719007
$ws.Range("F22").Value = 2.706230208981879
$ws.Range("G22").Value = 0.002497233902031499
$ws.Range("L22").Value = 0.2224103922807359
$ws.Range("M22").Value = 0.5350083953751721
$ws.Range("N22").Value = 2.212178830116486
$ws.Range("B23").Value = 3.048669247051237
$ws.Range("C23").Value = 0.4230485978395109
$ws.Range("D23").Value = 0.1096409244295486
$ws.Range("F23").Value = 2.675583033387397
$ws.Range("G23").Value = 0.002499538417779837
$ws.Range("L23").Value = 0.2213917092226794
$ws.Range("M23").Value = 0.5257528610446371
$ws.Range("N23").Value = 2.214991524128408
$ws.Range("B24").Value = 2.795587082890847
$ws.Range("C24").Value = 0.3626073556777669
$ws.Range("D24").Value = 0.1100249196773504
$ws.Range("F24").Value = 2.562273360431163
$ws.Range("G24").Value = 0.002508595985366978
$ws.Range("L24").Value = 0.2177085777268672
$ws.Range("M24").Value = 0.4912124000998688
$ws.Range("N24").Value = 2.226906965195454
$ws.Range("B25").Value = 2.528695678322151
$ws.Range("C25").Value = 0.2974852850019829
$ws.Range("D25").Value = 0.1106033064550616
$ws.Range("F25").Value = 2.44584491211242
$ws.Range("G25").Value = 0.002519076128250396
$ws.Range("L25").Value = 0.2141083523458533
$ws.Range("M25").Value = 0.4550362725868595
$ws.Range("N25").Value = 2.242425966404539
